$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the rich-text "CMS(µ)(2010)" string with a plain "CMS(mu)(2010)" string
# in all cells that reference it (L2:L7).
$ws.Range("L2:L7").Value = "CMS(mu)(2010)"

# Move the active selection from L13 to L16 to match the saved view state.
$ws.Range("L16").Select()
